# Add "LiPo charger" row to the off-board parts BOM (row 8, columns H:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(13).ColumnWidth = 16.7

$ws.Range("H8").Value = "LiPo charger"
$ws.Range("M8").Value = "Look for alternative"
$ws.Range("I8").Value = "adafruit.com/product/1944"
$ws.Range("J8").Value = 1
$ws.Range("L8").Value = 14.95

$ws.Hyperlinks.Add($ws.Range("I8"), "www.adafruit.com/product/259", [Type]::Missing, [Type]::Missing, "www.adafruit.com/product/259") | Out-Null

$ws.Range("L9").Select()
